$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.737.12"
$ws.Range("E2").Value = "'  +2.16%  "
$ws.Range("D3").Value = "'2.117.60"
$ws.Range("E3").Value = "'  +11.13%  "
$ws.Range("E4").Value = "'  +0.19%  "
$ws.Range("D5").Value = "'332.33"
$ws.Range("E5").Value = "'  +4.10%  "
$ws.Range("D6").Value = "'0.9999"
$ws.Range("E6").Value = "'  +0.08%  "
$ws.Range("D7").Value = "'0.5209"
$ws.Range("E7").Value = "'  +3.23%  "
$ws.Range("D8").Value = "'0.4368"
$ws.Range("E8").Value = "'  +7.10%  "
$ws.Range("D9").Value = "'0.09000"
$ws.Range("E9").Value = "'  +8.09%  "
$ws.Range("D10").Value = "'45.25"
$ws.Range("E10").Value = "'  +7.09%  "
$ws.Range("D11").Value = "'1.175"
$ws.Range("E11").Value = "'  +6.68%  "
$ws.Range("D12").Value = "'24.86"
$ws.Range("E12").Value = "'  +4.22%  "
$ws.Range("D13").Value = "'2.116.79"
$ws.Range("E13").Value = "'  +11.35%  "
$ws.Range("E14").Value = "'  +6.54%  "
$ws.Range("D15").Value = "'7.682"
$ws.Range("E15").Value = "'  +6.67%  "
$ws.Range("D16").Value = "'97.59"
$ws.Range("E16").Value = "'  +5.90%  "
$ws.Range("E17").Value = "'  +4.36%  "
$ws.Range("E18").Value = "'  -0.14%  "
$ws.Range("D19").Value = "'0.06619"
$ws.Range("E19").Value = "'  +2.12%  "
$ws.Range("D20").Value = "'19.24"
$ws.Range("E20").Value = "'  +5.39%  "
$ws.Range("D21").Value = "'6.453"
$ws.Range("E21").Value = "'  +9.04%  "
$ws.Range("D22").Value = "'0.9998"
$ws.Range("E22").Value = "'  +0.05%  "
$ws.Range("D23").Value = "'30.948.61"
$ws.Range("E23").Value = "'  +2.86%  "
$ws.Range("E24").Value = "'  +6.23%  "
$ws.Range("D25").Value = "'2.368.42"
$ws.Range("D26").Value = "'2.275"
$ws.Range("D27").Value = "'22.99"
$ws.Range("E27").Value = "'  +6.04%  "
$ws.Range("D28").Value = "'2.568"
$ws.Range("E28").Value = "'  +12.34%  "
$ws.Range("D29").Value = "'163.96"
$ws.Range("E29").Value = "'  +0.82%  "
$ws.Range("D30").Value = "'134.10"
$ws.Range("E30").Value = "'  +4.47%  "
$ws.Range("D31").Value = "'1.185"
$ws.Range("E31").Value = "'  +3.93%  "
$ws.Range("D32").Value = "'0.1070"
$ws.Range("E32").Value = "'  +2.89%  "
$ws.Range("D33").Value = "'6.237"
$ws.Range("E33").Value = "'  +4.91%  "
$ws.Range("B34").Value = "'HuobiToken"
$ws.Range("C34").Value = "'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "'3.906"
$ws.Range("E34").Value = "'  +3.64%  "
$ws.Range("B35").Value = "'ARBITRUM"
$ws.Range("C35").Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "'1.528"
$ws.Range("E35").Value = "'  +28.47%  "
$ws.Range("D36").Value = "'0.02586"
$ws.Range("E36").Value = "'  +5.64%  "
$ws.Range("E37").Value = "'  +5.12%  "
$ws.Range("D38").Value = "'0.06796"
$ws.Range("E38").Value = "'  +6.88%  "
$ws.Range("D39").Value = "'9.550"
$ws.Range("E39").Value = "'  +11.52%  "
$ws.Range("D40").Value = "'12.73"
$ws.Range("E40").Value = "'  +11.97%  "
$ws.Range("D41").Value = "'0.2253"
$ws.Range("E41").Value = "'  +5.32%  "
$ws.Range("E42").Value = "'  +4.34%  "
$ws.Range("D43").Value = "'1.257"
$ws.Range("E43").Value = "'  +3.98%  "
$ws.Range("D44").Value = "'14.45"
$ws.Range("E44").Value = "'  +8.08%  "
$ws.Range("D45").Value = "'0.9995"
$ws.Range("E45").Value = "'  +0.06%  "
$ws.Range("D46").Value = "'0.6320"
$ws.Range("E46").Value = "'  +4.43%  "
$ws.Range("D47").Value = "'2.254"
$ws.Range("E47").Value = "'  +2.50%  "
$ws.Range("D48").Value = "'3.667"
$ws.Range("E48").Value = "'  +1.62%  "
$ws.Range("D49").Value = "'1.278"
$ws.Range("E49").Value = "'  +5.83%  "
$ws.Range("D50").Value = "'127.38"
$ws.Range("E50").Value = "'  +4.99%  "
$ws.Range("D51").Value = "'83.32"
$ws.Range("E51").Value = "'  +5.72%  "
